$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells in D/E columns hold text-like values (prices with
# "." thousands separators, percentages with padding spaces) that must
# stay text, not get auto-converted to numbers by Excels type inference.
# Force text format before assigning so values round-trip as inlineStr/text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.985.38"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.051.02"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.71"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.610"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.69"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +6.75%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.386"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0812"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.24%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.351.65"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.59"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.81"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.30"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.754"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.045.87"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.902.39"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.28"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.86"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0837"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.79"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.32"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.55"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.133"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.04"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.33"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.53"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.59"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0612"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +9.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.32"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.08"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +9.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.30"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.48%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0217"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.482.73"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.06"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.71%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.57"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0922"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.30%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.09"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +13.98%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.97"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.05"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.242.28"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.70%  "
